$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$errorText = "Error: Message: stale element reference: stale element not found`r`n  (Session info: chrome=130.0.6723.70); For documentation on this error, please visit: https://www.selenium.dev/documentation/webdriver/troubleshooting/errors#stale-element-reference-exception`r`nStacktrace:`r`n`tGetHandleVerifier [0x00007FF616893AB5+28005]`r`n`t(No symbol) [0x00007FF6167F83B0]`r`n`t(No symbol) [0x00007FF61669580A]`r`n`t(No symbol) [0x00007FF6166A7D0B]`r`n`t(No symbol) [0x00007FF61669CF04]`r`n`t(No symbol) [0x00007FF61669AE09]`r`n`t(No symbol) [0x00007FF61669E5E7]`r`n`t(No symbol) [0x00007FF61669E6A0]`r`n`t(No symbol) [0x00007FF6166E692C]`r`n`t(No symbol) [0x00007FF6166E69F4]`r`n`t(No symbol) [0x00007FF6166DC808]`r`n`t(No symbol) [0x00007FF61670BA3A]`r`n`t(No symbol) [0x00007FF6166D9246]`r`n`t(No symbol) [0x00007FF61670BC50]`r`n`t(No symbol) [0x00007FF61672B8B3]`r`n`t(No symbol) [0x00007FF61670B7E3]`r`n`t(No symbol) [0x00007FF6166D75C8]`r`n`t(No symbol) [0x00007FF6166D8731]`r`n`tGetHandleVerifier [0x00007FF616B8643D+3118829]`r`n`tGetHandleVerifier [0x00007FF616BD6C90+3448640]`r`n`tGetHandleVerifier [0x00007FF616BCCF0D+3408317]`r`n`tGetHandleVerifier [0x00007FF61695A40B+841403]`r`n`t(No symbol) [0x00007FF61680340F]`r`n`t(No symbol) [0x00007FF6167FF484]`r`n`t(No symbol) [0x00007FF6167FF61D]`r`n`t(No symbol) [0x00007FF6167EEB79]`r`n`tBaseThreadInitThunk [0x00007FFB490D7374+20]`r`n`tRtlUserThreadStart [0x00007FFB4A61CC91+33]`r`n"

$ws.Range("B4").Value = $errorText
